$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "61.032.43"
$ws.Range("E2").Value2 = "  -1.09%  "
$ws.Range("D3").Value2 = "2.953.81"
$ws.Range("E3").Value2 = "  -1.44%  "
$ws.Range("E4").Value2 = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "587.06"
$ws.Range("E5").Value2 = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "147.54"
$ws.Range("E6").Value2 = "  +1.25%  "
$ws.Range("E7").Value2 = "  +0.11%  "
$ws.Range("D8").Value2 = "2.937.82"
$ws.Range("E8").Value2 = "  -1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.501"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "6.77"
$ws.Range("E10").Value2 = "  +11.79%  "
$ws.Range("E11").Value2 = "  -2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.450"
$ws.Range("E12").Value2 = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000225"
$ws.Range("E13").Value2 = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "34.41"
$ws.Range("E14").Value2 = "  +0.41%  "
$ws.Range("E15").Value2 = "  -0.71%  "
$ws.Range("D16").Value2 = "3.440.84"
$ws.Range("E16").Value2 = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "6.87"
$ws.Range("E17").Value2 = "  -0.80%  "
$ws.Range("D18").Value2 = "61.021.72"
$ws.Range("E18").Value2 = "  -1.12%  "
$ws.Range("D19").Value2 = "2.946.63"
$ws.Range("E19").Value2 = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "432.70"
$ws.Range("E20").Value2 = "  -3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.87"
$ws.Range("E21").Value2 = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.675"
$ws.Range("E22").Value2 = "  -1.35%  "
$ws.Range("E23").Value2 = "  -0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "80.34"
$ws.Range("E24").Value2 = "  -1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "11.03"
$ws.Range("E25").Value2 = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.21"
$ws.Range("E26").Value2 = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "11.89"
$ws.Range("E27").Value2 = "  -1.54%  "
$ws.Range("E28").Value2 = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.37"
$ws.Range("E29").Value2 = "  +1.89%  "
$ws.Range("B30").Value2 = "FirstDigitalUSD"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.00"
$ws.Range("E30").Value2 = "  -0.08%  "
$ws.Range("B31").Value2 = "PancakeSwap"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "2.65"
$ws.Range("E31").Value2 = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "2.19"
$ws.Range("E32").Value2 = "  +6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "26.88"
$ws.Range("E33").Value2 = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.107"
$ws.Range("E34").Value2 = "  -2.83%  "
$ws.Range("D35").Value2 = "0.0₃0838"
$ws.Range("E35").Value2 = "  +1.29%  "
$ws.Range("E36").Value2 = "  -1.14%  "
$ws.Range("E37").Value2 = "  -1.21%  "
$ws.Range("B38").Value2 = "OKB"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "50.01"
$ws.Range("E38").Value2 = "  -0.80%  "
$ws.Range("B39").Value2 = "dogwifhat"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.00"
$ws.Range("E39").Value2 = "  +4.79%  "
$ws.Range("B40").Value2 = "Stacks"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.03"
$ws.Range("E40").Value2 = "  +0.18%  "
$ws.Range("B41").Value2 = "Kaspa"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.126"
$ws.Range("E41").Value2 = "  +1.95%  "
$ws.Range("E42").Value2 = "  -2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.295"
$ws.Range("E43").Value2 = "  +9.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "42.54"
$ws.Range("E44").Value2 = "  +6.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0348"
$ws.Range("E45").Value2 = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "373.01"
$ws.Range("E46").Value2 = "  -5.74%  "
$ws.Range("D47").Value2 = "2.664.12"
$ws.Range("E47").Value2 = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "132.82"
$ws.Range("E48").Value2 = "  +0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "25.67"
$ws.Range("E49").Value2 = "  +9.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.15"
$ws.Range("E51").Value2 = "  -0.81%  "
